# Vietnamese translation pass for "Email T-1 [TEMPLATE] Partner email - if RSVP no"
$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($find, $replace) {
    # Whole-document, unique text replace (safe for runs that are not
    # immediately preceded by a zero-width marker such as a comment range
    # start or a hyperlink boundary).
    $r = $d.Content
    $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# 1. Brief paragraph
Replace-Text "An email sent to partners in the target country who have RSVPed no. It will be sent via customer.io" "Một email gửi tới các đối tác trong quốc gia mục tiêu đã không xác nhận tham gia. It will be sent via customer.io"

# 2. "Target audience" heading (bold)
Replace-Text "Target audience" "Đối tượng mục tiêu"

# 3. Target audience body
Replace-Text "Invited partners who RSVP no" "Các đối tác được mời đã không xác nhận tham gia"

# 4. "Subject line" label (bold)
Replace-Text "Subject line" "Tiêu đề"

# 5. Subject line text
Replace-Text ": Thinking of you at " ": Rất tiếc vì không thể gặp bạn tại sự kiện "

# 6. Delete the extra trailing plain-space run (duplicate highlighted/plain
#    space pair collapses into a single highlighted space run).
$r = $d.Content
$r.Find.Execute("[EVENT NAME]", $false) | Out-Null
$base = $r.End
$trail = $d.Range($base + 4, $base + 5)
if ($trail.Text -eq " ") {
    $trail.Text = ""
}

# 7. Headline "We'll miss you at the "
Replace-Text "We’ll miss you at the " "Chúng tôi rất tiếc vì bạn không thể tham dự sự kiện "

# 8. Remove the "Dear " run entirely
$r = $d.Content
$r.Find.Execute("Dear ", $false) | Out-Null
$r.Text = ""

# 9. ", " immediately after [PARTNER NAME] -> " thân mến, "
$r = $d.Content
$r.Find.Execute("[PARTNER NAME]", $false) | Out-Null
$tail = $d.Range($r.End, $r.End + 2)
if ($tail.Text -eq ", ") {
    $tail.Text = " thân mến, "
}

# 10. Thank-you paragraph lead-in
Replace-Text "Thank you for taking the time to respond to our invitation to the upcoming " "Cảm ơn bạn đã dành thời gian để phản hồi lời mời của chúng tôi cho sự kiện "

# 11. Thank-you paragraph tail
Replace-Text ". We were really looking forward to seeing you there." " sắp tới. Chúng tôi đã rất mong được gặp bạn tại sự kiện."

# 12. Disappointed paragraph
Replace-Text "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up. " "Mặc dù rất tiếc vì không thể gặp bạn, chúng tôi hiểu rằng bạn có công việc bận rộn và nhiều mối lo toan khác trong cuộc sống. "

# 13. Feedback paragraph
Replace-Text "If you’re comfortable sharing it with us, we’d like to know why you responded no. Please reply to this email as your feedback could help us make improvements in our event planning processes and better serve you in the future." "Nếu được, bạn có thể chia sẻ với chúng tôi lý do bạn không thể tham dự sự kiện không? Thông tin phản hồi của bạn sẽ giúp chúng tôi lên kế hoạch sự kiện tốt hơn cũng như hỗ trợ bạn hiệu quả hơn nữa trong tương lai."

# 14. Future events paragraph
Replace-Text "We hope to see you at our future events. " "Chúng tôi hy vọng sẽ được gặp bạn tại các sự kiện sắp tới của chúng tôi. "

# 15. "If you have any questions, please contact us via " - this run starts
#     right after <w:commentRangeStart/>, so a direct whole-range replace
#     would shift the marker past the new text. Fix the first character in
#     place, then replace the remainder, which keeps the marker anchored.
$r = $d.Content
$r.Find.Execute("If you have any questions, please contact us via ", $false) | Out-Null
$start = $r.Start
$newText = "Nếu bạn cần hỗ trợ, hãy liên hệ với chúng tôi qua "
$headRange = $d.Range($start, $start + 1)
$headRange.Text = $newText.Substring(0, 1)
$r2 = $d.Content
$r2.Start = $start + 1
$r2.Find.Execute("f you have any questions, please contact us via ", $false) | Out-Null
$r2.Text = $newText.Substring(1)

# 16. " or " between the live-chat and WhatsApp hyperlinks. The run starts
#     right after </w:hyperlink>, so skip the leading space (left untouched,
#     already correct) and replace starting from the "o" to avoid inheriting
#     hyperlink formatting.
$r = $d.Content
$r.Find.Execute("live chat", $false) | Out-Null
$base = $r.End
$orRange = $d.Range($base + 2, $base + 5)
if ($orRange.Text -eq "or ") {
    $orRange.Text = "hoặc "
}

# 17. "If you have any questions, please contact your country manager, "
Replace-Text "If you have any questions, please contact your country manager, " "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn "

# 18. ", at " after [NAME]
Replace-Text ", at " ", qua email "

# 19. " or " after [EMAIL ADDRESS]
$r = $d.Content
$r.Find.Execute("[EMAIL ADDRESS]", $false) | Out-Null
$base = $r.End
$orRange2 = $d.Range($base + 1, $base + 4)
if ($orRange2.Text -eq "or ") {
    $orRange2.Text = "hoặc số "
}
